# Updated cryptos list on Mon Feb 19 07:51:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (preserving things like leading/trailing
# zeros, thousand-separator dots, subscripted digits, % strings, etc.)
# without leaving a residual style index on the cell (so styles.xml /
# cell "s" attributes stay identical to the pre-edit file).
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

# --- Simple price / volume(1h) updates -------------------------------
Set-TextValue "D2" "52.394.62"
Set-TextValue "E2" "  +1.29%  "
Set-TextValue "D3" "2.921.64"
Set-TextValue "E3" "  +4.19%  "
Set-TextValue "D5" "353.04"
Set-TextValue "E5" "  -0.80%  "
Set-TextValue "D6" "112.98"
Set-TextValue "E6" "  +3.15%  "
Set-TextValue "D7" "0.562"
Set-TextValue "E7" "  +1.37%  "
Set-TextValue "E8" "  +0.04%  "
Set-TextValue "E9" "  +0.81%  "
Set-TextValue "D10" "40.36"
Set-TextValue "E10" "  +0.77%  "
Set-TextValue "E11" "  +3.10%  "
Set-TextValue "E12" "  +0.59%  "
Set-TextValue "D13" "20.19"
Set-TextValue "E13" "  +0.93%  "
Set-TextValue "D14" "7.89"
Set-TextValue "E14" "  +1.51%  "
Set-TextValue "D15" "3.386.23"
Set-TextValue "E15" "  +4.41%  "
Set-TextValue "D16" "2.939.02"
Set-TextValue "E16" "  +4.76%  "
Set-TextValue "D17" "0.996"
Set-TextValue "E17" "  +5.53%  "
Set-TextValue "D18" "52.431.16"
Set-TextValue "E18" "  +1.44%  "
Set-TextValue "D19" "7.68"
Set-TextValue "E19" "  -0.67%  "
Set-TextValue "D20" "3.33"
Set-TextValue "E20" "  +5.63%  "
Set-TextValue "D21" "14.47"
Set-TextValue "E21" "  +5.83%  "
Set-TextValue "D22" "0.0₃0986"
Set-TextValue "D23" "71.19"
Set-TextValue "E23" "  +1.11%  "
Set-TextValue "D24" "271.79"
Set-TextValue "E24" "  +1.24%  "
Set-TextValue "E25" "  +2.43%  "
Set-TextValue "D26" "27.10"
Set-TextValue "E26" "  +3.83%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.03%  "
Set-TextValue "E28" "  +0.45%  "
Set-TextValue "D29" "10.67"
Set-TextValue "E29" "  +2.78%  "
Set-TextValue "D30" "38.25"
Set-TextValue "E30" "  +2.72%  "
Set-TextValue "D31" "6.53"
Set-TextValue "E31" "  +5.12%  "
Set-TextValue "E32" "  +0.99%  "
Set-TextValue "D33" "6.19"
Set-TextValue "E33" "  +7.75%  "

# --- Rows 34 / 35 swap (OKB <-> Hedera) --------------------------------
Set-TextValue "B34" "Hedera"
Set-TextValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.0950"
Set-TextValue "E34" "  +10.58%  "

Set-TextValue "B35" "OKB"
Set-TextValue "C35" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D35" "53.10"
Set-TextValue "E35" "  +2.20%  "

Set-TextValue "D36" "0.0454"
Set-TextValue "E36" "  +1.52%  "
Set-TextValue "E37" "  -0.07%  "
Set-TextValue "D38" "3.35"
Set-TextValue "E38" "  +6.86%  "
Set-TextValue "D39" "19.00"
Set-TextValue "E39" "  +0.55%  "
Set-TextValue "E40" "  +4.18%  "
Set-TextValue "D41" "2.74"
Set-TextValue "E41" "  +10.17%  "
Set-TextValue "D42" "24.25"
Set-TextValue "E42" "  +10.92%  "
Set-TextValue "E43" "  +1.97%  "
Set-TextValue "D44" "122.59"
Set-TextValue "E44" "  +2.65%  "
Set-TextValue "D45" "2.60"
Set-TextValue "E46" "  +0.33%  "

# --- Rows 47 / 48 swap (NEARProtocol <-> Maker) ------------------------
Set-TextValue "B47" "Maker"
Set-TextValue "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "2.222.23"
Set-TextValue "E47" "  +4.55%  "

Set-TextValue "B48" "NEARProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D48" "3.57"
Set-TextValue "E48" "  +5.89%  "

Set-TextValue "D49" "0.266"
Set-TextValue "E49" "  +25.07%  "
Set-TextValue "D50" "0.0340"
Set-TextValue "E50" "  +15.13%  "
Set-TextValue "E51" "  +6.25%  "
